# Update "想去人数" (F column) figures to the freshly scraped values.
# Mapping derived from the OOXML diff: sheet name -> { row -> new value }

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ 3 = 1145; 7 = 801; 11 = 431; 14 = 951; 15 = 130; 16 = 2065; 18 = 9050; 21 = 82; 24 = 241 }
    "演出"     = @{ 12 = 63; 13 = 22 }
    "本地生活" = @{ 2 = 5619; 4 = 412 }
    "全部类型" = @{ 3 = 5619; 5 = 412; 7 = 1145; 12 = 801; 17 = 431; 22 = 951; 24 = 130; 27 = 2065; 29 = 9050; 30 = 63; 31 = 22; 34 = 82; 39 = 241 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
